$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.052715860460884
$ws.Range("D2").Value = 1.061967517465511
$ws.Range("E2").Value = 1.059912666932644
$ws.Range("F2").Value = 1.071623534516373
$ws.Range("I2").Value = 1.053178888479634
$ws.Range("J2").Value = 1.057737001665972
$ws.Range("K2").Value = 1.064690140997008
$ws.Range("L2").Value = 1.062640877421111
$ws.Range("M2").Value = 1.074320214451601
$ws.Range("N2").Value = 1.022993979595092
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.053758222802232
$ws.Range("D3").Value = 1.062673508311454
$ws.Range("E3").Value = 1.06080305019
$ws.Range("F3").Value = 1.072508628648978
$ws.Range("I3").Value = 1.05347952164289
$ws.Range("J3").Value = 1.058429507069345
$ws.Range("K3").Value = 1.06521112991984
$ws.Range("L3").Value = 1.063345390897089
$ws.Range("M3").Value = 1.075021728117258
$ws.Range("N3").Value = 1.023229410875736
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.054433076831815
$ws.Range("D4").Value = 1.063130545541598
$ws.Range("E4").Value = 1.061379823220677
$ws.Range("F4").Value = 1.073081929138189
$ws.Range("I4").Value = 1.053672926559786
$ws.Range("J4").Value = 1.058877373076295
$ws.Range("K4").Value = 1.065547782062059
$ws.Range("L4").Value = 1.063801255030026
$ws.Range("M4").Value = 1.07547559803923
$ws.Range("N4").Value = 1.023381552353683
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.054716875038016
$ws.Range("D5").Value = 1.063322733971634
$ws.Range("E5").Value = 1.06162244974679
$ws.Range("F5").Value = 1.073323083652609
$ws.Range("I5").Value = 1.053753964213377
$ws.Range("J5").Value = 1.059065599928741
$ws.Range("K5").Value = 1.065689199039169
$ws.Range("L5").Value = 1.063992898769816
$ws.Range("M5").Value = 1.075666390548616
$ws.Range("N5").Value = 1.023445464841545
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.054764531246152
$ws.Range("D6").Value = 1.063355006149383
$ws.Range("E6").Value = 1.061663196642194
$ws.Range("F6").Value = 1.073363582673444
$ws.Range("I6").Value = 1.053767554958323
$ws.Range("J6").Value = 1.059097200741973
$ws.Range("K6").Value = 1.065712937003782
$ws.Range("L6").Value = 1.064025076484051
$ws.Range("M6").Value = 1.075698424578962
$ws.Range("N6").Value = 1.023456193217935
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.054436868605886
$ws.Range("D7").Value = 1.063133113377985
$ws.Range("E7").Value = 1.061383064615555
$ws.Range("F7").Value = 1.073085150911933
$ws.Range("I7").Value = 1.053674010449506
$ws.Range("J7").Value = 1.05887988839362
$ws.Range("K7").Value = 1.065549672122146
$ws.Range("L7").Value = 1.063803815790085
$ws.Range("M7").Value = 1.075478147475837
$ws.Range("N7").Value = 1.023382406543648
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.053068053610719
$ws.Range("D8").Value = 1.062206065333576
$ws.Range("E8").Value = 1.060213443931006
$ws.Range("F8").Value = 1.071922534315543
$ws.Range("I8").Value = 1.053280721479436
$ws.Range("J8").Value = 1.057971084612779
$ws.Range("K8").Value = 1.064866306899378
$ws.Range("L8").Value = 1.062878971243229
$ws.Range("M8").Value = 1.074557305388998
$ws.Range("N8").Value = 1.023073585569136
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.050658923865025
$ws.Range("D9").Value = 1.060574181427516
$ws.Range("E9").Value = 1.058157334323927
$ws.Range("F9").Value = 1.069878389780482
$ws.Range("I9").Value = 1.052579102438478
$ws.Range("J9").Value = 1.056367913109092
$ws.Range("K9").Value = 1.063658625227009
$ws.Range("L9").Value = 1.061249286692129
$ws.Range("M9").Value = 1.072934273102435
$ws.Range("N9").Value = 1.02252789866036
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.049054807634413
$ws.Range("D10").Value = 1.05948747008068
$ws.Range("E10").Value = 1.056789954968985
$ws.Range("F10").Value = 1.068518743023322
$ws.Range("I10").Value = 1.052105606579908
$ws.Range("J10").Value = 1.055297996944893
$ws.Range("K10").Value = 1.062851197489703
$ws.Range("L10").Value = 1.060162881721428
$ws.Range("M10").Value = 1.071852042315926
$ws.Range("N10").Value = 1.022163113857794
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.048360678894503
$ws.Range("D11").Value = 1.059017214469353
$ws.Range("E11").Value = 1.05619867350257
$ws.Range("F11").Value = 1.067930753973172
$ws.Range("I11").Value = 1.051899219526611
$ws.Range("J11").Value = 1.054834449500794
$ws.Range("K11").Value = 1.062501035499141
$ws.Range("L11").Value = 1.05969247763122
$ws.Range("M11").Value = 1.071383385731499
$ws.Range("N11").Value = 1.022004926070335
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.048102918539144
$ws.Range("D12").Value = 1.058842586491686
$ws.Range("E12").Value = 1.055979166761391
$ws.Range("F12").Value = 1.067712461958684
$ws.Range("I12").Value = 1.051822354239159
$ws.Range("J12").Value = 1.054662227802604
$ws.Range("K12").Value = 1.062370889332669
$ws.Range("L12").Value = 1.059517751755222
$ws.Range("M12").Value = 1.071209300105334
$ws.Range("N12").Value = 1.021946133271536
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.048158205840588
$ws.Range("D13").Value = 1.058880042689815
$ws.Range("E13").Value = 1.056026246188781
$ws.Range("F13").Value = 1.067759281192518
$ws.Range("I13").Value = 1.051838851326622
$ws.Range("J13").Value = 1.054699171722156
$ws.Range("K13").Value = 1.062398809753393
$ws.Range("L13").Value = 1.059555230895377
$ws.Range("M13").Value = 1.071246642307812
$ws.Range("N13").Value = 1.021958746097259
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.048339370914522
$ws.Range("D14").Value = 1.05900277872315
$ws.Range("E14").Value = 1.056180526522598
$ws.Range("F14").Value = 1.067912707569674
$ws.Range("I14").Value = 1.051892869976031
$ws.Range("J14").Value = 1.054820214404995
$ws.Range("K14").Value = 1.062490279212268
$ws.Range("L14").Value = 1.059678034667656
$ws.Range("M14").Value = 1.071368995870867
$ws.Range("N14").Value = 1.022000066947772
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.048451001973701
$ws.Range("D15").Value = 1.059078406561044
$ws.Range("E15").Value = 1.056275599858211
$ws.Range("F15").Value = 1.068007253669732
$ws.Range("I15").Value = 1.051926125636232
$ws.Range("J15").Value = 1.054894787573486
$ws.Range("K15").Value = 1.062546625932829
$ws.Range("L15").Value = 1.059753698562258
$ws.Range("M15").Value = 1.07144438120838
$ws.Range("N15").Value = 1.02202552148305
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.049100884471746
$ws.Range("D16").Value = 1.059518685767787
$ws.Range("E16").Value = 1.0568292133645
$ws.Range("F16").Value = 1.06855778174964
$ws.Range("I16").Value = 1.052119275185679
$ws.Range("J16").Value = 1.055328755488405
$ws.Range("K16").Value = 1.062874425268051
$ws.Range("L16").Value = 1.060194101295788
$ws.Range("M16").Value = 1.071883144664322
$ws.Range("N16").Value = 1.022173607363187
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.049508663024782
$ws.Range("D17").Value = 1.059794941789885
$ws.Range("E17").Value = 1.057176696068565
$ws.Range("F17").Value = 1.068903314213207
$ws.Range("I17").Value = 1.052240068914942
$ws.Range("J17").Value = 1.055600901069424
$ws.Range("K17").Value = 1.063079900990211
$ws.Range("L17").Value = 1.060470359173471
$ws.Range("M17").Value = 1.072158358327008
$ws.Range("N17").Value = 1.022266435389183
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.049746558126138
$ws.Range("D18").Value = 1.059956105885652
$ws.Range("E18").Value = 1.057379454250233
$ws.Range("F18").Value = 1.069104929452725
$ws.Range("I18").Value = 1.052310394584151
$ws.Range("J18").Value = 1.055759613215739
$ws.Range("K18").Value = 1.063199699304162
$ws.Range("L18").Value = 1.060631497259366
$ws.Range("M18").Value = 1.072318881545981
$ws.Range("N18").Value = 1.022320557852292
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.049827681794929
$ws.Range("D19").Value = 1.060011063497774
$ws.Range("E19").Value = 1.057448602651096
$ws.Range("F19").Value = 1.069173687218835
$ws.Range("I19").Value = 1.052334351541917
$ws.Range("J19").Value = 1.055813725569614
$ws.Range("K19").Value = 1.063240538555729
$ws.Range("L19").Value = 1.060686441423
$ws.Range("M19").Value = 1.072373615068422
$ws.Range("N19").Value = 1.022339008379554
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.049464907612672
$ws.Range("D20").Value = 1.05976529916303
$ws.Range("E20").Value = 1.057139406450581
$ws.Range("F20").Value = 1.0688662343949
$ws.Range("I20").Value = 1.052227122465786
$ws.Range("J20").Value = 1.05557170508545
$ws.Range("K20").Value = 1.063057860783255
$ws.Range("L20").Value = 1.060440719153459
$ws.Range("M20").Value = 1.072128830939571
$ws.Range("N20").Value = 1.022256478151339
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.048286020382441
$ws.Range("D21").Value = 1.058966634744006
$ws.Range("E21").Value = 1.056135091447719
$ws.Range("F21").Value = 1.06786752418941
$ws.Range("I21").Value = 1.051876968456765
$ws.Range("J21").Value = 1.054784571438996
$ws.Range("K21").Value = 1.062463345949914
$ws.Range("L21").Value = 1.059641871925165
$ws.Range("M21").Value = 1.071332965946885
$ws.Range("N21").Value = 1.021987899945912
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.047545212102327
$ws.Range("D22").Value = 1.058464748842845
$ws.Range("E22").Value = 1.055504342284779
$ws.Range("F22").Value = 1.067240251489807
$ws.Range("I22").Value = 1.051655632967598
$ws.Range("J22").Value = 1.054289440782927
$ws.Range("K22").Value = 1.06208908602624
$ws.Range("L22").Value = 1.059139623516462
$ws.Range("M22").Value = 1.070832541451076
$ws.Range("N22").Value = 1.021818832899541
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.047937890129659
$ws.Range("D23").Value = 1.058730782434114
$ws.Range("E23").Value = 1.055838647359283
$ws.Range("F23").Value = 1.06757271807162
$ws.Range("I23").Value = 1.051773078788452
$ws.Range("J23").Value = 1.054551940393175
$ws.Range("K23").Value = 1.06228753209356
$ws.Range("L23").Value = 1.059405872915259
$ws.Range("M23").Value = 1.071097828755046
$ws.Range("N23").Value = 1.021908477524491
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.049484678661557
$ws.Range("D24").Value = 1.059778693301796
$ws.Range("E24").Value = 1.057156255784736
$ws.Range("F24").Value = 1.068882988946323
$ws.Range("I24").Value = 1.052232972814497
$ws.Range("J24").Value = 1.055584897573491
$ws.Range("K24").Value = 1.063067819965408
$ws.Range("L24").Value = 1.060454112197928
$ws.Range("M24").Value = 1.072142173108155
$ws.Range("N24").Value = 1.022260977468217
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.051281396034399
$ws.Range("D25").Value = 1.060995854367117
$ws.Range("E25").Value = 1.058688299913779
$ws.Range("F25").Value = 1.070406305823081
$ws.Range("I25").Value = 1.052761503284734
$ws.Range("J25").Value = 1.056782574153441
$ws.Range("K25").Value = 1.063971250093409
$ws.Range("L25").Value = 1.061670593442424
$ws.Range("M25").Value = 1.073353906657804
$ws.Range("N25").Value = 1.02266914826192
